$d = $word.ActiveDocument

$replacements = @(
    @("291÷9=32, 3", "453÷5=90, 3"),
    @("897÷6=149, 3", "297÷5=59, 2"),
    @("935÷5=187, 0", "559÷6=93, 1"),
    @("122÷3=40, 2", "682÷2=341, 0"),
    @("643÷8=80, 3", "216÷4=54, 0"),
    @("324÷8=40, 4", "359÷9=39, 8"),
    @("208÷6=34, 4", "719÷5=143, 4"),
    @("841÷3=280, 1", "411÷4=102, 3"),
    @("646÷2=323, 0", "237÷6=39, 3"),
    @("475÷6=79, 1", "925÷6=154, 1"),
    @("162÷7=23, 1", "656÷3=218, 2"),
    @("874÷9=97, 1", "293÷6=48, 5"),
    @("472÷8=59, 0", "761÷3=253, 2"),
    @("140÷2=70, 0", "817÷7=116, 5"),
    @("729÷4=182, 1", "344÷2=172, 0"),
    @("157÷3=52, 1", "776÷6=129, 2"),
    @("938÷2=469, 0", "778÷3=259, 1"),
    @("611÷3=203, 2", "163÷3=54, 1"),
    @("906÷3=302, 0", "405÷3=135, 0"),
    @("696÷5=139, 1", "262÷6=43, 4"),
    @("130÷8=16, 2", "955÷8=119, 3"),
    @("695÷4=173, 3", "803÷7=114, 5"),
    @("847÷5=169, 2", "955÷9=106, 1"),
    @("826÷8=103, 2", "912÷6=152, 0"),
    @("358÷8=44, 6", "522÷6=87, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
